# Make the CO column template placeholders ("CO{{ coXX }}") bold in the
# question-paper table, so the rendered CO values show up bold like the
# other header/columns do.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $cell = $t.Cell($r, 8)
    $txt = $cell.Range.Text
    if ($txt -like "*CO{{*") {
        $cell.Range.Font.Bold = 1
    }
}
